$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# --- 1. Update existing text: "scriptable spawn waves" -> "spawn wave with gaps" (row 26, col B) ---
$ws.Range("B26").Value2 = "spawn wave with gaps"

# --- 2. Insert a new row into the table (row 36) for "spawn point can be attacked" ---
# Shift rows 36:49 down to 37:50 (this also shifts the underlying table data rows).
$ws.Rows("36:36").Insert()

# Fill in the values for the newly inserted row 36.
$ws.Range("A36").Value2 = "Level"
$ws.Range("B36").Value2 = "spawn point"
$ws.Range("C36").Value2 = "spawn point can be attacked"
$ws.Range("D36").Value2 = "Code"

# --- 3. Resize the worksheet table ("Tabelle1") to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A6:D50"))

# --- 4. Cosmetic: adjust column A width (best-fit recalculated by Excel) ---
$ws.Columns.Item(1).ColumnWidth = 58.3

# --- 5. Cosmetic: update selection / active cell shown when the workbook was saved ---
$ws.Range("D37").Select()
